$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LoginData")

# The DRM service URL now includes an explicit scheme.
$ws.Range("C2").Value = "https://qanintendo.dsidrm.com"

# Turn the cell into a clickable hyperlink (shortened display text/address).
$ws.Hyperlinks.Add($ws.Range("C2"), "https://qanintendo") | Out-Null
$link = $ws.Hyperlinks.Item(1)
$link.TextToDisplay = "https://qanintendo"

# Re-assert the full URL as the cell's text - TextToDisplay above only
# updates the hyperlink's cached display text, not the stored value.
$ws.Range("C2").Value = "https://qanintendo.dsidrm.com"

# The row grew slightly taller once the link was added.
$ws.Rows.Item(2).RowHeight = 14.9

"done"
